# Update "想去人数" (interest count, column F) values across the four
# sheets of the workbook to match the refreshed data snapshot.
# Commit: "Update gh-pages to output generated at 456a3b4"
#
# NOTE: this COM-interop engine does not propagate *named* parameters
# into functions correctly, so all helper calls below use positional
# arguments.

function Set-ColumnFValues {
    param(
        $SheetName,
        $RowUpdates
    )

    $ws = $wb.Worksheets.Item($SheetName)
    foreach ($row in $RowUpdates.Keys) {
        $ws.Cells.Item($row, 6).Value = $RowUpdates[$row]
    }
}

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions)
Set-ColumnFValues "展览" @{
    2  = 1175
    3  = 1087
    4  = 1889
    6  = 1239
    7  = 67
    8  = 25
    9  = 133
    10 = 328
    11 = 108
    13 = 784
    14 = 224
    15 = 120
    16 = 32
    18 = 341
    19 = 205
    20 = 689
    21 = 66
    24 = 46
    25 = 899
    26 = 343
    27 = 186
    29 = 299
}

# 演出 (Performances)
Set-ColumnFValues "演出" @{
    4  = 329
    5  = 20
    10 = 622
}

# 本地生活 (Local life)
Set-ColumnFValues "本地生活" @{
    2 = 324
}

# 全部类型 (All types)
Set-ColumnFValues "全部类型" @{
    2  = 324
    3  = 1175
    4  = 1087
    5  = 1889
    7  = 1239
    8  = 67
    10 = 25
    11 = 133
    12 = 328
    13 = 108
    15 = 784
    16 = 224
    17 = 120
    19 = 32
    20 = 329
    22 = 20
    23 = 341
    27 = 205
    28 = 689
    29 = 66
    32 = 46
    33 = 899
    34 = 343
    37 = 186
    39 = 299
    40 = 622
}
